$d = $word.ActiveDocument

$pairs = @(
    @("51÷8=6, 3", "85÷5=17, 0"),
    @("11÷2=5, 1", "62÷8=7, 6"),
    @("68÷9=7, 5", "36÷9=4, 0"),
    @("74÷4=18, 2", "82÷5=16, 2"),
    @("51÷7=7, 2", "76÷4=19, 0"),
    @("96÷3=32, 0", "51÷8=6, 3"),
    @("27÷4=6, 3", "47÷6=7, 5"),
    @("12÷7=1, 5", "89÷2=44, 1"),
    @("81÷2=40, 1", "80÷4=20, 0"),
    @("29÷6=4, 5", "17÷8=2, 1"),
    @("87÷6=14, 3", "33÷7=4, 5"),
    @("26÷4=6, 2", "92÷6=15, 2"),
    @("67÷2=33, 1", "35÷6=5, 5"),
    @("12÷4=3, 0", "28÷7=4, 0"),
    @("83÷9=9, 2", "85÷2=42, 1"),
    @("42÷4=10, 2", "96÷2=48, 0"),
    @("47÷9=5, 2", "68÷2=34, 0"),
    @("21÷2=10, 1", "17÷7=2, 3"),
    @("12÷3=4, 0", "54÷6=9, 0"),
    @("60÷6=10, 0", "97÷7=13, 6"),
    @("35÷9=3, 8", "84÷8=10, 4"),
    @("10÷7=1, 3", "64÷3=21, 1"),
    @("28÷8=3, 4", "36÷3=12, 0"),
    @("75÷7=10, 5", "58÷2=29, 0"),
    @("30÷7=4, 2", "13÷4=3, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
